# Discharge workbook update: add a "New depth" sensitivity block to STN2 and
# STN3, re-using the already-computed "D" (depth) and "V" (velocity) columns
# from the earlier blocks on each sheet, then leave the workbook scrolled to
# STN3 (matching the final author view state).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# STN2 ("sheet2"): append a new "New depth" block at rows 43-60
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("STN2")

$ws2.Range("A43").Value = "New depth"
$ws2.Range("A43").Font.Bold = $true

$ws2.Range("A44").Value = "X"
$ws2.Range("B44").Value = "D"
$ws2.Range("C44").Value = "V"
$ws2.Range("D44").Value = "segment"
$ws2.Range("E44").Value = "Q"
$ws2.Range("F44").Value = "Qtotal"

$ws2.Range("A45").Value = 0.55
$ws2.Range("B45").Formula = "=B23*2.54"
$ws2.Range("C45").Value = 0
$ws2.Range("D45").Formula = "=A45"
$ws2.Range("F45").Formula = "=SUM(E45:E63)"

$ws2.Range("A46").Value = 0.6
$ws2.Range("B46").Formula = "=B24*2.54"
$ws2.Range("C46").Value = 0
$ws2.Range("D46").Formula = "=(A46+(A47-A46)/2)"
$ws2.Range("E46").Formula = "=(D46-D45)*(B46)*C46"

$ws2.Range("A47").Value = 0.63
$ws2.Range("B47").Formula = "=B25*2.54"
$ws2.Range("C47").Value = 0
$ws2.Range("D47").Formula = "=(A47+(A48-A47)/2)"
$ws2.Range("E47").Formula = "=(D47-D46)*(B47)*C47"

$ws2.Range("A48").Value = 0.66
$ws2.Range("B48").Formula = "=B26*2.54"
$ws2.Range("C48").Value = 0
$ws2.Range("D48").Formula = "=(A48+(A49-A48)/2)"
$ws2.Range("E48").Formula = "=(D48-D47)*(B48)*C48"

$ws2.Range("A49").Value = 0.69
$ws2.Range("B49").Formula = "=B27*2.54"
$ws2.Range("C49").Value = 0
$ws2.Range("D49").Formula = "=(A49+(A50-A49)/2)"
$ws2.Range("E49").Formula = "=(D49-D48)*(B49)*C49"

$ws2.Range("A50").Value = 0.72
$ws2.Range("B50").Formula = "=B28*2.54"
$ws2.Range("C50").Value = 0.0572
$ws2.Range("D50").Formula = "=(A50+(A51-A50)/2)"
$ws2.Range("E50").Formula = "=(D50-D49)*(B50)*C50"

$ws2.Range("A51").Value = 0.75
$ws2.Range("B51").Formula = "=B29*2.54"
$ws2.Range("C51").Value = 0.0572
$ws2.Range("D51").Formula = "=(A51+(A52-A51)/2)"
$ws2.Range("E51").Formula = "=(D51-D50)*(B51)*C51"

$ws2.Range("A52").Value = 0.78
$ws2.Range("B52").Formula = "=B30*2.54"
$ws2.Range("C52").Value = 0.0572
$ws2.Range("D52").Formula = "=(A52+(A53-A52)/2)"
$ws2.Range("E52").Formula = "=(D52-D51)*(B52)*C52"

$ws2.Range("A53").Value = 0.81
$ws2.Range("B53").Formula = "=B31*2.54"
$ws2.Range("C53").Value = 0.06864
$ws2.Range("D53").Formula = "=(A53+(A54-A53)/2)"
$ws2.Range("E53").Formula = "=(D53-D52)*(B53)*C53"

$ws2.Range("A54").Value = 0.84
$ws2.Range("B54").Formula = "=B32*2.54"
$ws2.Range("C54").Value = 0.0858
$ws2.Range("D54").Formula = "=(A54+(A55-A54)/2)"
$ws2.Range("E54").Formula = "=(D54-D53)*(B54)*C54"

$ws2.Range("A55").Value = 0.87
$ws2.Range("B55").Formula = "=B33*2.54"
$ws2.Range("C55").Value = 0.10296
$ws2.Range("D55").Formula = "=(A55+(A56-A55)/2)"
$ws2.Range("E55").Formula = "=(D55-D54)*(B55)*C55"

$ws2.Range("A56").Value = 0.9
$ws2.Range("B56").Formula = "=B34*2.54"
$ws2.Range("C56").Value = 0.10868
$ws2.Range("D56").Formula = "=(A56+(A57-A56)/2)"
$ws2.Range("E56").Formula = "=(D56-D55)*(B56)*C56"

$ws2.Range("A57").Value = 0.93
$ws2.Range("B57").Formula = "=B35*2.54"
$ws2.Range("C57").Value = 0.10868
$ws2.Range("D57").Formula = "=(A57+(A58-A57)/2)"
$ws2.Range("E57").Formula = "=(D57-D56)*(B57)*C57"

$ws2.Range("A58").Value = 0.96
$ws2.Range("B58").Formula = "=B36*2.54"
$ws2.Range("C58").Value = 0.09724
$ws2.Range("D58").Formula = "=(A58+(A59-A58)/2)"
$ws2.Range("E58").Formula = "=(D58-D57)*(B58)*C58"

$ws2.Range("A59").Value = 0.99
$ws2.Range("B59").Formula = "=B37*2.54"
$ws2.Range("C59").Value = 0.07436
$ws2.Range("D59").Formula = "=(A59+(A60-A59)/2)"
$ws2.Range("E59").Formula = "=(D59-D58)*(B59)*C59"

$ws2.Range("A60").Value = 1.03
$ws2.Range("B60").Formula = "=B38*2.54"
$ws2.Range("C60").Value = 0
$ws2.Range("D60").Formula = "=(A60+(A61-A60)/2)"
$ws2.Range("E60").Formula = "=(D60-D59)*(B60)*C60"

# ---------------------------------------------------------------------------
# STN3 ("sheet3"): append a new "new depth" block at rows 31-44
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("STN3")

$ws3.Range("A31").Value = "new depth"
$ws3.Range("A31").Font.Bold = $true

$ws3.Range("A32").Value = "X"
$ws3.Range("B32").Value = "D"
$ws3.Range("C32").Value = "V"
$ws3.Range("D32").Value = "segment"
$ws3.Range("E32").Value = "Q"
$ws3.Range("F32").Value = "Qtotal"

$ws3.Range("A33").Value = 0.34
$ws3.Range("B33").Formula = "=B17*2.54"
$ws3.Range("C33").Value = 0
$ws3.Range("D33").Formula = "=A33"
$ws3.Range("F33").Formula = "=SUM(E33:E51)"

$ws3.Range("A34").Value = 0.36
$ws3.Range("B34").Formula = "=B18*2.54"
$ws3.Range("C34").Value = 0.15444000000000002
$ws3.Range("D34").Formula = "=(A34+(A35-A34)/2)"
$ws3.Range("E34").Formula = "=(D34-D33)*(B34)*C34"

$ws3.Range("A35").Value = 0.39
$ws3.Range("B35").Formula = "=B19*2.54"
$ws3.Range("C35").Value = 0.429
$ws3.Range("D35").Formula = "=(A35+(A36-A35)/2)"
$ws3.Range("E35").Formula = "=(D35-D34)*(B35)*C35"

$ws3.Range("A36").Value = 0.42
$ws3.Range("B36").Formula = "=B20*2.54"
$ws3.Range("C36").Value = 0.58916000000000002
$ws3.Range("D36").Formula = "=(A36+(A37-A36)/2)"
$ws3.Range("E36").Formula = "=(D36-D35)*(B36)*C36"

$ws3.Range("A37").Value = 0.45
$ws3.Range("B37").Formula = "=B21*2.54"
$ws3.Range("C37").Value = 0.5434
$ws3.Range("D37").Formula = "=(A37+(A38-A37)/2)"
$ws3.Range("E37").Formula = "=(D37-D36)*(B37)*C37"

$ws3.Range("A38").Value = 0.48
$ws3.Range("B38").Formula = "=B22*2.54"
$ws3.Range("C38").Value = 0.51480000000000004
$ws3.Range("D38").Formula = "=(A38+(A39-A38)/2)"
$ws3.Range("E38").Formula = "=(D38-D37)*(B38)*C38"

$ws3.Range("A39").Value = 0.51
$ws3.Range("B39").Formula = "=B23*2.54"
$ws3.Range("C39").Value = 0.38324000000000003
$ws3.Range("D39").Formula = "=(A39+(A40-A39)/2)"
$ws3.Range("E39").Formula = "=(D39-D38)*(B39)*C39"

$ws3.Range("A40").Value = 0.54
$ws3.Range("B40").Formula = "=B24*2.54"
$ws3.Range("C40").Value = 0.35464000000000001
$ws3.Range("D40").Formula = "=(A40+(A41-A40)/2)"
$ws3.Range("E40").Formula = "=(D40-D39)*(B40)*C40"

$ws3.Range("A41").Value = 0.57
$ws3.Range("B41").Formula = "=B25*2.54"
$ws3.Range("C41").Value = 0.27455999999999997
$ws3.Range("D41").Formula = "=(A41+(A42-A41)/2)"
$ws3.Range("E41").Formula = "=(D41-D40)*(B41)*C41"

$ws3.Range("A42").Value = 0.6
$ws3.Range("B42").Formula = "=B26*2.54"
$ws3.Range("C42").Value = 0
$ws3.Range("D42").Formula = "=(A42+(A43-A42)/2)"
$ws3.Range("E42").Formula = "=(D42-D41)*(B42)*C42"

$ws3.Range("A43").Value = 0.62
$ws3.Range("B43").Formula = "=B27*2.54"
$ws3.Range("C43").Value = 0
$ws3.Range("D43").Formula = "=(A43+(A44-A43)/2)"
$ws3.Range("E43").Formula = "=(D43-D42)*(B43)*C43"

$ws3.Range("D44").Formula = "=(A44+(A45-A44)/2)"
$ws3.Range("E44").Formula = "=(D44-D43)*(B44)*C44"

# ---------------------------------------------------------------------------
# View state: STN1 scrolled down, STN2 selection on its new total cell, and
# STN3 becomes the active/selected tab (matches activeTab="2" in workbook.xml)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("STN1")
$ws1.Range("A18").Select()

$ws2.Range("F45").Select()

$ws3.Range("F33").Select()
$ws3.Activate()
